# Update the "Metadata" sheet: Version, Status, Date and Contact values.
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value  = "0.4.0-snapshot-1"
$meta.Range("B6").Value  = "draft"
$meta.Range("B8").Value  = "2024-05-23T12:16:26+00:00"
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"

# Update the "Elements" sheet: the two rightmost "Mapping" columns (AK, AL)
# swap places - the "Mapping: RIM Mapping" column and the
# "Mapping: Spécification métier..." column exchange their header text and
# their per-row values.
$els = $wb.Worksheets.Item("Elements")

$akHeader = $els.Cells.Item(1, 37)
$alHeader = $els.Cells.Item(1, 38)
$akHeaderVal = $akHeader.Value2
$alHeaderVal = $alHeader.Value2
$akHeader.Value = $alHeaderVal
$alHeader.Value = $akHeaderVal

for ($r = 2; $r -le 24; $r++) {
    $akCell = $els.Cells.Item($r, 37)
    $alCell = $els.Cells.Item($r, 38)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

# The two columns' widths swap along with their content (column AK becomes
# the wide column, AL becomes the narrow one).
$els.Columns.Item(37).ColumnWidth = 81.1
$els.Columns.Item(38).ColumnWidth = 24.1
